$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Member 3 effort entries for days 4-7 (rows 11-14) -- the actual data edit.
# Everything else on the sheet (B1 total, B3 average, the E/F/G formula
# columns and the chart's cached series) is formula-driven off these cells
# and recalculates automatically.
$ws.Range("D11").Value = 1
$ws.Range("D12").Value = 2
$ws.Range("D13").Value = 5
$ws.Range("D14").Value = 5

# Selection moved to D14 as last-active cell.
$ws.Range("D14").Select()
